# Daily scrape update - 2026-02-22 04:23:18 UTC
# Refresh the "global-talent" opportunity listing: rows 2-8 get new
# opportunity data, the old rows 9-11 are dropped, a couple of column
# widths are retuned, and the one-off "Yes"/yellow-highlight premium
# flag on E2 goes back to the plain "No" look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Drop the trailing rows (9, 10, 11) that no longer exist in today's
#    scrape. This also shrinks the sheet's used range to A1:H8.
# ---------------------------------------------------------------------
$ws.Range("A9:H11").EntireRow.Delete()

# ---------------------------------------------------------------------
# 2) Re-tune a handful of column widths.
#    (COM ColumnWidth is expressed in "characters" and Excel stores the
#    width in the xlsx as ColumnWidth + 5/6, so back that offset out to
#    land exactly on the target stored widths of 39 / 29 / 15 / 34.)
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 39 - 5/6
$ws.Columns.Item(4).ColumnWidth = 29 - 5/6
$ws.Columns.Item(6).ColumnWidth = 15 - 5/6
$ws.Columns.Item(8).ColumnWidth = 34 - 5/6

# ---------------------------------------------------------------------
# 3) Helper: write a value that must stay plain text even though it
#    looks numeric (opportunity IDs) - a leading apostrophe forces text
#    entry, then resetting the style back to Normal strips the
#    resulting quote-prefix formatting so no stray style sticks around.
# ---------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 4) Row 2
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A2") "1332079"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1332079"
$ws.Range("C2").Value = "Language Specialist - Portuguese"
$ws.Range("D2").Value = "Colombo, Sri Lanka"
$ws.Range("E2").Value = "No"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "3 - 6 Months"
$ws.Range("H2").Value = "Aitken Spence Travels (Pvt) Ltd"

# ---------------------------------------------------------------------
# 5) Row 3
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A3") "1331934"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1331934"
$ws.Range("C3").Value = "Content Creator and Marketing Intern"
$ws.Range("D3").Value = "Lahore, Pakistan"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "Plush Natural"

# ---------------------------------------------------------------------
# 6) Row 4
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A4") "1331452"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1331452"
$ws.Range("C4").Value = "Digital Marketing Specialist"
$ws.Range("D4").Value = "Mumbai, Maharashtra, India"
$ws.Range("F4").Value = "0 applicants"
$ws.Range("H4").Value = "Supervity AI"

# ---------------------------------------------------------------------
# 7) Row 5
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A5") "1331450"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1331450"
$ws.Range("C5").Value = "Academy & Community Lead"
$ws.Range("D5").Value = "Mumbai, Maharashtra, India"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("H5").Value = "Supervity AI"

# ---------------------------------------------------------------------
# 8) Row 6
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A6") "1331440"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1331440"
$ws.Range("C6").Value = "Forward Deployed Engineer"
$ws.Range("D6").Value = "Mumbai, Maharashtra, India"
$ws.Range("F6").Value = "0 applicants"
$ws.Range("H6").Value = "Supervity AI"

# ---------------------------------------------------------------------
# 9) Row 7
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A7") "1331174"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1331174"
$ws.Range("C7").Value = "Graphic Designer / UI - UX Designer"
$ws.Range("D7").Value = "Karachi, Pakistan"
$ws.Range("F7").Value = "4 applicants"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "Wavetec"

# ---------------------------------------------------------------------
# 10) Row 8
# ---------------------------------------------------------------------
Set-TextValue $ws.Range("A8") "1326614"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1326614"
$ws.Range("C8").Value = "Language Specialist – Italian"
$ws.Range("D8").Value = "Colombo, Sri Lanka"
$ws.Range("F8").Value = "0 applicants"
$ws.Range("G8").Value = "3 - 6 Months"
$ws.Range("H8").Value = "Aitken Spence Travels (Pvt) Ltd"
